# modify some data and folder name
# Re-write the "sounds-meta-data" sheet: a new row (sneezing / not-need now)
# is inserted after row 3, the rows that follow shift down by one, and a
# trailing row 15 (num=13, otherwise blank) is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sounds-meta-data")

# --- Row 4 : new "sneezing" entry ----------------------------------------------------
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 21
$ws.Range("C4").Value = "sneezing"
$ws.Range("D4").Value = "not-need now"
$ws.Range("E4").ClearContents()
$ws.Range("E4").Style = "Normal"

# --- Row 5 : coughing -------------------------------------------------------------------
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 24
$ws.Range("C5").Value = "coughing"
$ws.Range("D5").Clear()
$ws.Range("E5").Value = "*"

# --- Row 6 : kettle-sound ----------------------------------------------------------------
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 51
$ws.Range("C6").Value = "kettle-sound"
$ws.Range("D6").Value = "熱水壺笛聲"
$ws.Range("E6").Value = "*"

# --- Row 7 : alarm -----------------------------------------------------------------------
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 52
$ws.Range("C7").Value = "alarm"
$ws.Range("D7").Value = "All alert sounds"
$ws.Range("E7").Value = "*"

# --- Row 8 : boiling-water-bubble-sound ---------------------------------------------------
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 53
$ws.Range("C8").Value = "boiling-water-bubble-sound"
$ws.Range("D8").Clear()

# --- Row 9 : ringtone --------------------------------------------------------------------
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 54
$ws.Range("C9").Value = "ringtone"
$ws.Range("D9").Value = "Line and WeChat default"

# --- Row 10 : shower-water ----------------------------------------------------------------
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 55
$ws.Range("C10").Value = "shower-water"
$ws.Range("D10").Clear()

# --- Row 11 : pain-sound -------------------------------------------------------------------
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 57
$ws.Range("C11").Value = "pain-sound"
$ws.Range("D11").Style = "Normal"

# --- Row 12 : Foot-step --------------------------------------------------------------------
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 58
$ws.Range("C12").Value = "Foot-step"

# --- Row 13 : silence ----------------------------------------------------------------------
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 98
$ws.Range("C13").Value = "silence"

# --- Row 14 : other-sounds -------------------------------------------------------------------
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 99
$ws.Range("C14").Value = "other-sounds"

# --- Row 15 : new trailing (mostly blank) row -------------------------------------------------
$ws.Range("A15").Value = 13
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Style = "Normal"

# --- Row heights ---------------------------------------------------------------------------
$ws.Range("A4").EntireRow.RowHeight = 22.05
$ws.Range("A5").EntireRow.RowHeight = 22.05
$ws.Range("A11").EntireRow.RowHeight = 22.05
$ws.Range("A12").EntireRow.RowHeight = 22.05
$ws.Range("A13").EntireRow.RowHeight = 22.05
$ws.Range("A14").EntireRow.RowHeight = 22.05
$ws.Range("A15").EntireRow.RowHeight = 22.05

# --- Selection -------------------------------------------------------------------------------
$ws.Range("B18").Select() | Out-Null
